$d = $word.ActiveDocument

$replacements = @(
    @("2024-12-25 Wednesday", "2024-12-26 Thursday"),
    @("34×13=", "87×62="),
    @("57×15=", "25×71="),
    @("51×30=", "65×84="),
    @("17×62=", "26×38="),
    @("29×18=", "69×57="),
    @("36×32=", "31×56="),
    @("21×69=", "39×44="),
    @("55×70=", "40×97="),
    @("51×42=", "67×58="),
    @("27×25=", "24×75="),
    @("49×84=", "19×33="),
    @("18×95=", "65×85="),
    @("21×79=", "30×71="),
    @("74×92=", "86×90="),
    @("38×61=", "90×96="),
    @("60×53=", "34×36="),
    @("77×23=", "50×61="),
    @("24×47=", "21×98="),
    @("24×37=", "73×89="),
    @("56×82=", "15×23="),
    @("37×11=", "67×52="),
    @("81×77=", "42×92="),
    @("24×54=", "75×38="),
    @("66×66=", "56×95="),
    @("77×93=", "30×92=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
